$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @{ B=1.199213271610233;  C=0.354944670868008;   D=0.01503616514069961; E=0.4254321993331303; F=0.4030827437275448; I=0.2470587248881824; O=1.244608197349265 }
  3  = @{ B=1.049321093537515;  C=0.3110170542605033;  D=0.01326176328310424; E=0.3710241705299779; F=0.3993247100541879; I=0.2546024141071861; O=1.246152473051836 }
  4  = @{ B=0.9569891392220029; C=0.2839267843903315;  D=0.01216786834112327; E=0.3377147263827851; F=0.3975013686066049; I=0.2595819383436453; O=1.248714745250382 }
  5  = @{ B=0.9192903332054243; C=0.2728579909621658;  D=0.01172102449260137; E=0.3241633335824616; F=0.3968793584726313; I=0.2616981020691878; O=1.25016221876777  }
  6  = @{ B=0.9130261397797312; C=0.2710182727439019;  D=0.01164676265852904; E=0.3219144347038565; F=0.3967833622882466; I=0.2620547286666621; O=1.250426856040676 }
  7  = @{ B=0.9564810119189815; C=0.2837776246635428;  D=0.0121618463455988;  E=0.3375318791631088; F=0.3974924909217421; I=0.2596101261489689; O=1.248732636715403 }
  8  = @{ B=1.14759337227639;   C=0.3398233259204062;  D=0.0144252858462508;  E=0.4066508972364602; F=0.4016860639633606; I=0.2495873564182904; O=1.244804124213232 }
  9  = @{ B=1.519932936973134;  C=0.4487710321283203;  D=0.01882763574549529; E=0.5430615204112854; F=0.4137842399471268; I=0.2327115326399998; O=1.250024152912715 }
  10 = @{ B=1.791942598005846;  C=0.5282161006501269;  D=0.02203855149714684; E=0.6439546468780861; F=0.4250827099878194; I=0.2220342732705323; O=1.26190253960533  }
  11 = @{ B=1.915338907632815;  C=0.564225176901175;   D=0.02349391636342801; E=0.690031129686048;  F=0.4307563817382345; I=0.2175565732668723; O=1.269088469333695 }
  12 = @{ B=1.962015122589662;  C=0.5778416822854524;  D=0.02404423550894563; E=0.7075073357755883; F=0.4329824448454502; I=0.2159160425858815; O=1.272068743956197 }
  13 = @{ B=1.951964878313788;  C=0.574909993150527;   D=0.02392575037173827; E=0.7037422395795971; F=0.4324995599823183; I=0.2162669024471313; O=1.271415319235615 }
  14 = @{ B=1.91918002760093;   C=0.5653458057975058;  D=0.02353920759617267; E=0.6914683329632823; F=0.4309379622260394; I=0.2174204990424968; O=1.269328449053091 }
  15 = @{ B=1.899091619540798;  C=0.5594849293671018;  D=0.02330233432530093; E=0.6839539284820262; F=0.4299915633601756; I=0.2181342982195957; O=1.268084009353515 }
  16 = @{ B=1.783871293892219;  C=0.5258601436306094;  D=0.02194333044912611; E=0.6409472439179638; F=0.4247227336731996; I=0.2223345807100117; O=1.261469028813877 }
  17 = @{ B=1.713098210472424;  C=0.5051985243355261;  D=0.02110824433166414; E=0.6146116184445276; F=0.4216277917657152; I=0.2250088668850321; O=1.257869465458867 }
  18 = @{ B=1.672359312328695;  C=0.4933022189781013;  D=0.0206274283238983;  E=0.5994808435457628; F=0.4198978730123741; I=0.2265827293310387; O=1.255966690492045 }
  19 = @{ B=1.658560373303601;  C=0.4892722401429523;  D=0.02046454805777387; E=0.5943606385196745; F=0.4193207519033493; I=0.2271217251800426; O=1.255351143183873 }
  20 = @{ B=1.720635462592952;  C=0.5073992646434249;  D=0.02119719231695427; E=0.617413338128884;  F=0.4219520513896455; I=0.2247204877097619; O=1.258235278979754 }
  21 = @{ B=1.928811141306539;  C=0.5681555687934292;  D=0.02365276645058856; E=0.695072698124676;  F=0.4313945299252424; I=0.2170801608760939; O=1.269934358543964 }
  22 = @{ B=2.064566090058918;  C=0.6077502989650725;  D=0.02525297134520343; E=0.7459920709235348; F=0.4380182199325873; I=0.2124080368787808; O=1.279092037482741 }
  23 = @{ B=1.992139235232116;  C=0.5866283605292324;  D=0.02439934891990703; E=0.7187996576595594; F=0.434441369289047;  I=0.2148720707879548; O=1.27406516419245  }
  24 = @{ B=1.717228028083412;  C=0.5064043649381915;  D=0.02115698115353126; E=0.6161466499250565; F=0.4218052997839123; I=0.2248507506545945; O=1.258069375765729 }
  25 = @{ B=1.419472465225738;  C=0.4194019072779156;  D=0.01764071084648577; E=0.506051184457263;  F=0.4100916148856655; I=0.2369765749567776; O=1.247211425952941 }
}

foreach ($r in $data.Keys) {
  $row = $data[$r]
  $ws.Range("B$r").Value = $row.B
  $ws.Range("C$r").Value = $row.C
  $ws.Range("D$r").Value = $row.D
  $ws.Range("E$r").Value = $row.E
  $ws.Range("F$r").Value = $row.F
  $ws.Range("I$r").Value = $row.I
  $ws.Range("O$r").Value = $row.O
}
